$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text shown for the "R10" rule row (cell E8)
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection on the sheet as it was when saved
$ws.Range("E8").Select()
